$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell (C column, probabilities) without
# Excel's automatic "looks like a number" coercion turning it into a
# numeric cell. We compute the literal text via a throw-away formula cell
# far outside the used range, then PasteSpecial just the value into the
# real destination (keeping the destination's own pre-existing style).
function Set-TextValue($cell, [string]$text) {
    $helper = $ws.Cells.Item(5000, 1)
    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $cell.PasteSpecial(-4163) | Out-Null
    $helper.EntireRow.Delete()
}

# --- Batch 1: first five new spies (rows 8-12) -----------------------
# Entered column by column (Nombre, then Compañero, then Probabilidad),
# matching how the shared-strings table was built in the source edit.
$batch1 = @(
    @(8,  "Jose",     "Mabel",    "0.9"),
    @(9,  "Alvaro",   "Pepe",     "0.5"),
    @(10, "Julieta",  "Ivan",     "0.1"),
    @(11, "Mabel",    "Gabriela", "0.9"),
    @(12, "Gabriela", "Alvaro",   "0.7")
)

foreach ($entry in $batch1) { $ws.Cells.Item($entry[0], 1).Value = $entry[1] }
foreach ($entry in $batch1) { $ws.Cells.Item($entry[0], 2).Value = $entry[2] }
foreach ($entry in $batch1) { Set-TextValue $ws.Cells.Item($entry[0], 3) $entry[3] }

# --- Batch 2: three more spies (rows 13-15) ---------------------------
$batch2 = @(
    @(13, "Hugo",    "Ruben", "0.4"),
    @(14, "Candela", "Mabel", "0.8"),
    @(15, "Ruben",   "Juan",  "0.9")
)

foreach ($entry in $batch2) { $ws.Cells.Item($entry[0], 1).Value = $entry[1] }
foreach ($entry in $batch2) { $ws.Cells.Item($entry[0], 2).Value = $entry[2] }
foreach ($entry in $batch2) { Set-TextValue $ws.Cells.Item($entry[0], 3) $entry[3] }

$ws.Range("C15").Select()
